$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data cleaning: B21/B22 had the Indonesia/India country-name counterparts
# swapped relative to column A. Fix the lookup so each row's imf name
# counterpart matches its own country.
$ws.Range("B21").Value = "Indonesia"
$ws.Range("B22").Value = "India"

# Reflect the author's last selection in the sheet before saving.
$ws.Range("I20").Select()
